$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Shannon Diversity column values were being exponentiated; restore the
# natural (non-exponentiated) estimates and confidence intervals for the
# normal (Shannon Diversity) model rows.

# (Intercept) row
$t.Cell(3, 8).Range.Text = "3.13"
$t.Cell(3, 9).Range.Text = "(1.15, 5.1)"

# Number of Native Flowering Species row
$t.Cell(4, 8).Range.Text = "0.16"
$t.Cell(4, 9).Range.Text = "(-0.11, 0.44)"

# Average Native Bloom Cover row
$t.Cell(5, 8).Range.Text = "-0.12"
$t.Cell(5, 9).Range.Text = "(-0.37, 0.14)"
